# Definicion de significado de los estados
$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Fix typo "Delegate,stop" -> "delegate,stop" on slide 1
#    (matches the lower-case convention used by sibling shapes).
# ------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$fixShape = $null
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $candidate = $s1.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.TextRange.Text -eq "Delegate,stop") {
        $fixShape = $candidate
        break
    }
}
if ($fixShape -ne $null) {
    $firstLetter = $fixShape.TextFrame.TextRange.Characters(1, 1)
    $firstLetter.Text = "d"
}

# ------------------------------------------------------------------
# 2) Add a new slide 2: "Definición de los estados" with a 3x2 table
#    describing the READY / RESERVED states.
# ------------------------------------------------------------------
$s2 = $p.Slides.Add($p.Slides.Count + 1, 2)

$title = $s2.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Definición de los estados"
$title.TextFrame.TextRange.LanguageID = "es-ES"

# Remove the empty body placeholder and replace it with a table,
# mirroring the "insert table into content placeholder" gesture.
$s2.Shapes.Item(2).Delete()

$tblShape = $s2.Shapes.AddTable(3, 2, 66, 143.75, 828, 130)
$tblShape.Name = "Marcador de contenido 3"

$table = $tblShape.Table
$table.Columns.Item(1).Width = 145.5
$table.Columns.Item(2).Width = 682.5
$table.Rows.Item(1).Height = 29.2
$table.Rows.Item(2).Height = 29.2
$table.Rows.Item(3).Height = 29.2

function Set-CellText($row, $col, $text) {
    $tr = $table.Cell($row, $col).Shape.TextFrame.TextRange
    $tr.Text = $text
    $tr.LanguageID = "es-ES"
}

Set-CellText 1 1 "ESTADO"
Set-CellText 1 2 "DESCRIPTION"
Set-CellText 2 1 "READY"
Set-CellText 2 2 "La tarea se encuentra lista para que un miembro del grupo de participantes que ha sido nominado para realizar la tarea pueda reclamarla y con ello reservarla para ejecutarla."
Set-CellText 3 1 "RESERVED"
Set-CellText 3 2 "La tarea ha sido asignada a un participante específico para su ejecución. Los demás participante ya no podrán ejecutar la tarea."
